# Update the "as_of_utc" timestamp column (AA) for the data rows (2-26)
# on both the "Главные" and "Линейные" sheets from
# "2025-11-19 03:11:26" to "2025-11-19 07:08:30".

$wb = $excel.ActiveWorkbook

$sheetNames = @("Главные", "Линейные")
$newTimestamp = "2025-11-19 07:08:30"

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    for ($row = 2; $row -le 26; $row++) {
        $ws.Cells.Item($row, 27).Value = $newTimestamp
    }
}
